# "foundation done for vba replica"
#
# Populates the "turf" sheet (gameweek tracker for the FPL replica) with
# the two currently-tracked managers' full line-up/captain picks
# (columns T:AH) and backfills the rest of the league table (rows 4:33,
# columns A:B) that already exists on the "FPL" sheet.

$wb  = $excel.ActiveWorkbook
$fpl  = $wb.Worksheets.Item("FPL")
$turf = $wb.Worksheets.Item("turf")

# ---------------------------------------------------------------------
# Row 2 - "MarxUnited"
# ---------------------------------------------------------------------
$turf.Range("A2").Value = "MarxUnited"
$turf.Range("B2").Value = 2438499

$turf.Range("T2").Value  = "Sánchez 6"
$turf.Range("U2").Value  = "Duffy 6"
$turf.Range("V2").Value  = "Livramento 4"
$turf.Range("W2").Value  = "White 7"
$turf.Range("X2").Value  = "Salah 13"
$turf.Range("Y2").Value  = "Gray 2"
$turf.Range("Z2").Value  = "Raphinha 3"
$turf.Range("AA2").Value = "Sarr 1"
$turf.Range("AB2").Value = "Antonio 2"
$turf.Range("AC2").Value = "Ronaldo 1"
$turf.Range("AD2").Value = "Lukaku 4"
$turf.Range("AE2").Value = "Steer 0"
$turf.Range("AF2").Value = "S.Longstaff 2"
$turf.Range("AG2").Value = "Alexander-Arnold 0"
$turf.Range("AH2").Value = "Amartey 0"

# ---------------------------------------------------------------------
# Row 3 - "Lucky"
# ---------------------------------------------------------------------
$turf.Range("A3").Value = "Lucky"
$turf.Range("B3").Value = 1293900

$turf.Range("T3").Value  = "Ramsdale 6"
$turf.Range("U3").Value  = "Dias 0"
$turf.Range("V3").Value  = "Rüdiger 2"
$turf.Range("W3").Value  = "Cancelo 0"
$turf.Range("X3").Value  = "Saka 3"
$turf.Range("Y3").Value  = "Gallagher 2"
$turf.Range("Z3").Value  = "Salah 13"
$turf.Range("AA3").Value = "Raphinha 3"
$turf.Range("AB3").Value = "Lukaku 4"
$turf.Range("AC3").Value = "Saint-Maximin 2"
$turf.Range("AD3").Value = "Antonio 2"
$turf.Range("AE3").Value = "Foster 2"
$turf.Range("AF3").Value = "Christensen 0"
$turf.Range("AG3").Value = "Douglas Luiz 2"
$turf.Range("AH3").Value = "Alexander-Arnold 0"

# ---------------------------------------------------------------------
# Rows 4-33 - remainder of the league table, mirrored from "FPL"
# (PasteSpecial values-only so text like "False"/team tags round-trips
# as text instead of being re-typed and coerced to bool/number).
# ---------------------------------------------------------------------
$fpl.Range("A4:B33").Copy()
$turf.Range("A4:B33").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Match "FPL"'s Manager_iD number styling (style: Consolas 9, blue) on
# the same rows, without touching A2:B3 which were typed in plain.
$fpl.Range("B2").Copy()
$turf.Range("B4:B33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# View state: selection now spans the populated table, A2:B33.
# ---------------------------------------------------------------------
$turf.Activate()
$turf.Range("A2:B33").Select()
